$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "CAD file" column (C) for every Straight Beams row (3-16) already holds
# "iam v1.0". Bump it to "iam v2.0" and, while we're at it, fill in the two
# rows that were still blank (S4x4 / S8x4, rows 17-18) with the same note -
# straight beams were redesigned as a modular buildup, so SXx4/SXx2/SXx1 are
# all now "ok" against the new CAD rev.
$ws.Range("C3:C18").Value = "iam v2.0"

# Move the selection / top-left-cell to match the saved view state.
$ws.Range("C26").Select()
